$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Update J column (2025) values for OPM (Adultas Desaparecidas row2, Adultas Localizadas row3)
# and CAVV (Menores Desaparecidas row5, Menores Localizadas row6)
$ws.Range("J2").Value = 160
$ws.Range("J3").Value = 152
$ws.Range("J5").Value = 115
$ws.Range("J6").Value = 108

# Update the selected cell/range to match the author's last selection
$ws.Range("K12").Select()
